# Apply the cryptos-list refresh described by the commit diff: update
# Price (D) / Volume(1h) (E) figures, and fix the two swapped coin rows
# (B/C/D/E for rows 38-39 and 47-48).
#
# Numeric-looking Price strings (e.g. "532.21") are written with a leading
# apostrophe so Excel keeps them as literal text - matching the source
# file's inlineStr cells - instead of silently converting them to a
# Number; the apostrophe itself is not part of the stored value. The
# Style is then reset to "Normal" on those cells so no stray
# quote-prefix/number-format style is left behind (the source cells carry
# no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.203.37"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "3.189.31"
$ws.Range("E3").Value = "  -4.84%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'532.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.07%  "
$ws.Range("D6").Value = "'134.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.72%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "3.187.41"
$ws.Range("E8").Value = "  -4.93%  "
$ws.Range("D9").Value = "'0.454"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.12%  "
$ws.Range("D10").Value = "'7.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.06%  "
$ws.Range("E11").Value = "  -7.80%  "
$ws.Range("E12").Value = "  -4.64%  "
$ws.Range("D13").Value = "3.734.01"
$ws.Range("E13").Value = "  -4.79%  "
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "'25.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.44%  "
$ws.Range("D16").Value = "3.189.46"
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("D17").Value = "58.337.24"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("E18").Value = "  -8.14%  "
$ws.Range("D19").Value = "'5.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.02%  "
$ws.Range("D20").Value = "'13.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.40%  "
$ws.Range("D21").Value = "'8.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.61%  "
$ws.Range("D22").Value = "'359.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.38%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'69.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.78%  "
$ws.Range("E25").Value = "  -7.84%  "
$ws.Range("D26").Value = "3.321.95"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "'0.168"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.50%  "
$ws.Range("D28").Value = "0.0₃0950"
$ws.Range("E28").Value = "  -12.10%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "'6.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.63%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'1.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.55%  "
$ws.Range("D33").Value = "'6.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.73%  "
$ws.Range("D34").Value = "'21.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("E35").Value = "  -7.70%  "
$ws.Range("D36").Value = "'4.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.28%  "
$ws.Range("D37").Value = "'160.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'6.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.67%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.75%  "
$ws.Range("D40").Value = "'25.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.68%  "
$ws.Range("D41").Value = "'0.0704"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.39%  "
$ws.Range("D42").Value = "3.215.93"
$ws.Range("E42").Value = "  -5.02%  "
$ws.Range("D43").Value = "'40.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").Value = "'0.706"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("E45").Value = "  -6.81%  "
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.06%  "
$ws.Range("D49").Value = "2.280.14"
$ws.Range("E49").Value = "  -7.60%  "
$ws.Range("D50").Value = "'6.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.49%  "
$ws.Range("D51").Value = "'20.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.62%  "
